$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = 67488109
$ws.Range("C2").Value = "Behöver inte valideras"
$ws.Range("P2").Value = "Koghult, 900 m SV , Sk"
$ws.Range("S2").Value = 50
$ws.Range("X2").Value = "M-Sbo-0024"
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2017-06-30"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2017-06-30"
$ws.Range("AY2").Value = "Floraväkteri Sverige"
